# Update the crypto price/volume/hour table to reflect the latest GitHub Actions scrape.
# All data cells in the sheet are stored as text, so values are written back as text
# (prefixed with a leading apostrophe) to avoid Excel coercing numeric-looking
# strings (e.g. "0.1420", "16") into actual numbers and losing formatting such as
# trailing zeros.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Cell = "D2"; Value = '245.60' }
    @{ Cell = "G2"; Value = '16' }
    @{ Cell = "G3"; Value = '16' }
    @{ Cell = "D4"; Value = '5.393' }
    @{ Cell = "G4"; Value = '16' }
    @{ Cell = "D5"; Value = '0.05856' }
    @{ Cell = "G5"; Value = '16' }
    @{ Cell = "D6"; Value = '3.395' }
    @{ Cell = "G6"; Value = '16' }
    @{ Cell = "D7"; Value = '6.384' }
    @{ Cell = "G7"; Value = '16' }
    @{ Cell = "D8"; Value = '0.8185' }
    @{ Cell = "G8"; Value = '16' }
    @{ Cell = "D9"; Value = '0.9982' }
    @{ Cell = "G9"; Value = '16' }
    @{ Cell = "B10"; Value = 'One' }
    @{ Cell = "C10"; Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one' }
    @{ Cell = "D10"; Value = '0.0005891' }
    @{ Cell = "E10"; Value = '9OneONE' }
    @{ Cell = "G10"; Value = '16' }
    @{ Cell = "B11"; Value = 'WazirX' }
    @{ Cell = "C11"; Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx' }
    @{ Cell = "D11"; Value = '0.1420' }
    @{ Cell = "E11"; Value = '10WazirXWRX' }
    @{ Cell = "G11"; Value = '16' }
    @{ Cell = "B12"; Value = 'LiechtensteinCryptoassetsExchange' }
    @{ Cell = "C12"; Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx' }
    @{ Cell = "D12"; Value = '0.03819' }
    @{ Cell = "E12"; Value = '11LiechtensteinCryptoassetsExchangeLCX' }
    @{ Cell = "G12"; Value = '16' }
    @{ Cell = "B13"; Value = 'MandalaExchangeToken' }
    @{ Cell = "C13"; Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx' }
    @{ Cell = "D13"; Value = '0.07403' }
    @{ Cell = "E13"; Value = '12MandalaExchangeTokenMDX' }
    @{ Cell = "G13"; Value = '16' }
    @{ Cell = "B14"; Value = 'BitrueCoin' }
    @{ Cell = "C14"; Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr' }
    @{ Cell = "D14"; Value = '0.03042' }
    @{ Cell = "E14"; Value = '13BitrueCoinBTR' }
    @{ Cell = "G14"; Value = '16' }
    @{ Cell = "B15"; Value = 'MCDex' }
    @{ Cell = "C15"; Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb' }
    @{ Cell = "D15"; Value = '4.135' }
    @{ Cell = "E15"; Value = '14MCDexMCB' }
    @{ Cell = "G15"; Value = '16' }
    @{ Cell = "B16"; Value = 'BitMartToken' }
    @{ Cell = "C16"; Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx' }
    @{ Cell = "D16"; Value = '0.09399' }
    @{ Cell = "E16"; Value = '15BitMartTokenBMX' }
    @{ Cell = "G16"; Value = '16' }
    @{ Cell = "B17"; Value = 'BitForexToken' }
    @{ Cell = "C17"; Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf' }
    @{ Cell = "D17"; Value = '0.001598' }
    @{ Cell = "E17"; Value = '16BitForexTokenBF' }
    @{ Cell = "G17"; Value = '16' }
    @{ Cell = "B18"; Value = 'CoinExToken' }
    @{ Cell = "C18"; Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet' }
    @{ Cell = "D18"; Value = '0.04826' }
    @{ Cell = "E18"; Value = '17CoinExTokenCET' }
    @{ Cell = "G18"; Value = '16' }
    @{ Cell = "D19"; Value = '0.005960' }
    @{ Cell = "G19"; Value = '16' }
    @{ Cell = "D20"; Value = '0.004084' }
    @{ Cell = "G20"; Value = '16' }
    @{ Cell = "D21"; Value = '0.0009938' }
    @{ Cell = "G21"; Value = '16' }
    @{ Cell = "D22"; Value = '0.0001500' }
    @{ Cell = "G22"; Value = '16' }
    @{ Cell = "D23"; Value = '3.740' }
    @{ Cell = "G23"; Value = '16' }
    @{ Cell = "D24"; Value = '2.221' }
    @{ Cell = "G24"; Value = '16' }
    @{ Cell = "D25"; Value = '0.3238' }
    @{ Cell = "G25"; Value = '16' }
    @{ Cell = "G26"; Value = '16' }
    @{ Cell = "D27"; Value = '0.0002494' }
    @{ Cell = "E27"; Value = '26UpBotsUBXTWorstin24h' }
    @{ Cell = "G27"; Value = '16' }
    @{ Cell = "G28"; Value = '16' }
    @{ Cell = "G29"; Value = '16' }
    @{ Cell = "G30"; Value = '16' }
    @{ Cell = "G31"; Value = '16' }
    @{ Cell = "G32"; Value = '16' }
    @{ Cell = "G33"; Value = '16' }
    @{ Cell = "G34"; Value = '16' }
    @{ Cell = "G35"; Value = '16' }
    @{ Cell = "G36"; Value = '16' }
    @{ Cell = "G37"; Value = '16' }
    @{ Cell = "G38"; Value = '16' }
    @{ Cell = "G39"; Value = '16' }
    @{ Cell = "D40"; Value = '0.03867' }
    @{ Cell = "G40"; Value = '16' }
    @{ Cell = "D41"; Value = '0.006411' }
    @{ Cell = "G41"; Value = '16' }
    @{ Cell = "G42"; Value = '16' }
    @{ Cell = "D43"; Value = '0.003001' }
    @{ Cell = "E43"; Value = '42CEJICEJIBestin24h' }
    @{ Cell = "G43"; Value = '16' }
    @{ Cell = "D44"; Value = '0.006672' }
    @{ Cell = "G44"; Value = '16' }
    @{ Cell = "D45"; Value = '0.00005624' }
    @{ Cell = "G45"; Value = '16' }
    @{ Cell = "G46"; Value = '16' }
    @{ Cell = "D47"; Value = '0.6501' }
    @{ Cell = "G47"; Value = '16' }
    @{ Cell = "D48"; Value = '0.1425' }
    @{ Cell = "G48"; Value = '16' }
    @{ Cell = "D49"; Value = '0.00002100' }
    @{ Cell = "G49"; Value = '16' }
    @{ Cell = "G50"; Value = '16' }
    @{ Cell = "G51"; Value = '16' }
)

foreach ($change in $changes) {
    $ws.Range($change.Cell).Value = "'" + $change.Value
}

